$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.6
$ws.Range("E2").Value = 0.738
$ws.Range("F2").Value = 0.737
$ws.Range("G2").Value = 0.29
$ws.Range("H2").Value = 1.482
$ws.Range("I2").Value = 3

# Row 3
$ws.Range("D3").Value = 0.6
$ws.Range("E3").Value = 1.016
$ws.Range("F3").Value = 0.935
$ws.Range("G3").Value = 0.264
$ws.Range("H3").Value = 1.788
$ws.Range("I3").Value = 3

# Row 4
$ws.Range("C4").Value = 8
$ws.Range("E4").Value = 1.012
$ws.Range("F4").Value = 1.065
$ws.Range("G4").Value = 0.216
$ws.Range("H4").Value = 1.585
$ws.Range("I4").Value = 3

# Row 5
$ws.Range("D5").Value = 0.628
$ws.Range("E5").Value = 0.449
$ws.Range("F5").Value = 0.443
$ws.Range("G5").Value = 0.609
$ws.Range("H5").Value = 1.516
$ws.Range("I5").Value = 3

# Row 6
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 0.766
$ws.Range("E6").Value = 1.47
$ws.Range("F6").Value = 1.611
$ws.Range("G6").Value = -0.077
$ws.Range("H6").Value = 2.278
$ws.Range("I6").Value = 3
